$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H29").Value = 5012.7144
$ws.Range("J29").Value = 7104
$ws.Range("L29").Value = 21312
$ws.Range("N29").Value = -21874

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 116.86667
$ws.Range("I33").Value = 82.42856999999999
$ws.Range("K33").Value = 82.42856999999999
$ws.Range("M33").Value = 146.57143

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 1644.3077
$ws.Range("I38").Value = 940.5454999999999
$ws.Range("J38").Value = 5515
$ws.Range("K38").Value = 2821.6365
$ws.Range("L38").Value = 16545
$ws.Range("M38").Value = -2449.6365
$ws.Range("N38").Value = -17289

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H55").Value = 1073.5
$ws.Range("I55").Value = 473
$ws.Range("J55").Value = 1974.25
$ws.Range("K55").Value = 473
$ws.Range("L55").Value = 1974.25
$ws.Range("M55").Value = -259
$ws.Range("N55").Value = -2402.25

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H58").Value = 1496.8334
$ws.Range("J58").Value = 2143.5
$ws.Range("L58").Value = 6430.5
$ws.Range("N58").Value = -6730.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H111").Value = 2457.75
$ws.Range("I111").Value = 994
$ws.Range("J111").Value = 3921.5
$ws.Range("K111").Value = 2982
$ws.Range("L111").Value = 11764.5
$ws.Range("M111").Value = 85
$ws.Range("N111").Value = -17898.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 6326.625
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 6326.625
$ws.Range("K138").Value = 0
$ws.Range("L138").Value = 18979.875
$ws.Range("M138").ClearContents()
$ws.Range("N138").Value = -29259.875

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H44").Value = 11772
$ws.Range("J44").Value = 11772
$ws.Range("L44").Value = 11772
$ws.Range("N44").Value = -12748

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2711.4119
$ws.Range("I45").Value = 2006.3334
$ws.Range("K45").Value = 2006.3334
$ws.Range("M45").Value = -1629.3334

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 3914.5
$ws.Range("I61").Value = 3528.3076
$ws.Range("K61").Value = 3528.3076
$ws.Range("M61").Value = -3316.3076

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 2497.7646
$ws.Range("I122").Value = 2579.1875
$ws.Range("J122").Value = 1195
$ws.Range("K122").Value = 7737.5625
$ws.Range("L122").Value = 3585
$ws.Range("M122").Value = -5287.5625
$ws.Range("N122").Value = -8485

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 3914.5
$ws.Range("I136").Value = 3528.3076
$ws.Range("K136").Value = 10584.9228
$ws.Range("M136").Value = -8034.9228

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H14").Value = 547.5
$ws.Range("J14").Value = 547.5
$ws.Range("L14").Value = 547.5
$ws.Range("N14").Value = -891.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 9999.5
$ws.Range("I86").Value = 0
$ws.Range("J86").Value = 9999.5
$ws.Range("K86").Value = 0
$ws.Range("L86").Value = 9999.5
$ws.Range("M86").ClearContents()
$ws.Range("N86").Value = -12245.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 9999.5
$ws.Range("I89").Value = 0
$ws.Range("J89").Value = 9999.5
$ws.Range("K89").Value = 0
$ws.Range("L89").Value = 49997.5
$ws.Range("M89").ClearContents()
$ws.Range("N89").Value = -61229.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H12").Value = 322.33334
$ws.Range("J12").Value = 314.42856
$ws.Range("L12").Value = 314.42856
$ws.Range("N12").Value = -654.4285600000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H47").Value = 0
$ws.Range("I47").Value = 0
$ws.Range("K47").Value = 0
$ws.Range("M47").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1806.3334
$ws.Range("I58").Value = 968.2
$ws.Range("K58").Value = 968.2
$ws.Range("M58").Value = -765.2

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H70").Value = 0
$ws.Range("J70").Value = 0
$ws.Range("L70").Value = 0
$ws.Range("N70").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H73").Value = 0
$ws.Range("J73").Value = 0
$ws.Range("L73").Value = 0
$ws.Range("N73").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 534.9
$ws.Range("J107").Value = 669.25
$ws.Range("L107").Value = 669.25
$ws.Range("N107").Value = -4509.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 1806.3334
$ws.Range("I136").Value = 968.2
$ws.Range("K136").Value = 2904.6
$ws.Range("M136").Value = -354.6000000000004

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H138").Value = 104998.336
$ws.Range("J138").Value = 104998.336
$ws.Range("L138").Value = 104998.336
$ws.Range("N138").Value = -115278.336

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H55").Value = 4052.6667
$ws.Range("I55").Value = 551
$ws.Range("J55").Value = 5803.5
$ws.Range("K55").Value = 1653
$ws.Range("L55").Value = 17410.5
$ws.Range("M55").Value = -1476
$ws.Range("N55").Value = -17764.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H70").Value = 2506
$ws.Range("J70").Value = 4000
$ws.Range("L70").Value = 12000
$ws.Range("N70").Value = -12630

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H73").Value = 2506
$ws.Range("J73").Value = 4000
$ws.Range("L73").Value = 12000
$ws.Range("N73").Value = -14184

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H74").Value = 0
$ws.Range("J74").Value = 0
$ws.Range("L74").Value = 0
$ws.Range("N74").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H77").Value = 0
$ws.Range("J77").Value = 0
$ws.Range("L77").Value = 0
$ws.Range("N77").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 773.1429000000001
$ws.Range("I16").Value = 888.2
$ws.Range("J16").Value = 485.5
$ws.Range("K16").Value = 888.2
$ws.Range("L16").Value = 485.5
$ws.Range("M16").Value = -718.2
$ws.Range("N16").Value = -825.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H18").Value = 9411.706
$ws.Range("I18").Value = 10461.538
$ws.Range("J18").Value = 5999.75
$ws.Range("K18").Value = 10461.538
$ws.Range("L18").Value = 5999.75
$ws.Range("M18").Value = -10289.538
$ws.Range("N18").Value = -6343.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H20").Value = 2583.3333
$ws.Range("J20").Value = 2583.3333
$ws.Range("L20").Value = 2583.3333
$ws.Range("N20").Value = -3035.3333

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2569.6
$ws.Range("I22").Value = 300
$ws.Range("J22").Value = 4082.6667
$ws.Range("K22").Value = 300
$ws.Range("L22").Value = 4082.6667
$ws.Range("M22").Value = -5
$ws.Range("N22").Value = -4672.6667

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H27").Value = 2569.6
$ws.Range("I27").Value = 300
$ws.Range("J27").Value = 4082.6667
$ws.Range("K27").Value = 300
$ws.Range("L27").Value = 4082.6667
$ws.Range("M27").Value = -193
$ws.Range("N27").Value = -4296.6667

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 4561.12
$ws.Range("I46").Value = 3190.7693
$ws.Range("K46").Value = 3190.7693
$ws.Range("M46").Value = -3002.7693

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H62").Value = 5000
$ws.Range("J62").Value = 5000
$ws.Range("L62").Value = 5000
$ws.Range("N62").Value = -6248

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H65").Value = 5000
$ws.Range("J65").Value = 5000
$ws.Range("L65").Value = 15000
$ws.Range("N65").Value = -21240

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H63").Value = 13800
$ws.Range("J63").Value = 13800
$ws.Range("L63").Value = 13800
$ws.Range("N63").Value = -15048

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H66").Value = 13800
$ws.Range("J66").Value = 13800
$ws.Range("L66").Value = 41400
$ws.Range("N66").Value = -47640

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H95").Value = 10000
$ws.Range("J95").Value = 10000
$ws.Range("L95").Value = 10000
$ws.Range("N95").Value = -15492

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H97").Value = 37500
$ws.Range("J97").Value = 37500
$ws.Range("L97").Value = 37500
$ws.Range("N97").Value = -39482
